$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.326.04"
$ws.Range("E2").Value = "  -2.63%  "
$ws.Range("D3").Value = "1.568.66"
$ws.Range("E3").Value = "  -3.53%  "
$ws.Range("E4").Value = "  -0.57%  "
$ws.Range("D5").Value = "'207.70"
$ws.Range("E5").Value = "  -2.92%  "
$ws.Range("E6").Value = "  -0.54%  "
$ws.Range("D7").Value = "'0.477"
$ws.Range("E7").Value = "  -5.03%  "
$ws.Range("E8").Value = "  -1.95%  "
$ws.Range("E9").Value = "  -1.54%  "
$ws.Range("D10").Value = "'17.89"
$ws.Range("E10").Value = "  -1.35%  "
$ws.Range("E11").Value = "  -0.85%  "
$ws.Range("D12").Value = "1.785.43"
$ws.Range("E12").Value = "  -3.69%  "
$ws.Range("D13").Value = "1.570.00"
$ws.Range("E13").Value = "  -3.66%  "
$ws.Range("E14").Value = "  -2.98%  "
$ws.Range("D15").Value = "'0.507"
$ws.Range("E15").Value = "  -2.72%  "
$ws.Range("D16").Value = "25.306.02"
$ws.Range("E16").Value = "  -2.63%  "
$ws.Range("D17").Value = "'59.65"
$ws.Range("E17").Value = "  -2.46%  "
$ws.Range("D18").Value = "0.0₃0710"
$ws.Range("E18").Value = "  -3.78%  "
$ws.Range("D20").Value = "'185.54"
$ws.Range("E20").Value = "  -2.03%  "
$ws.Range("E21").Value = "  -1.86%  "
$ws.Range("D22").Value = "'9.32"
$ws.Range("E22").Value = "  -2.28%  "
$ws.Range("E23").Value = "  -2.62%  "
$ws.Range("D24").Value = "'0.130"
$ws.Range("E24").Value = "  -1.63%  "
$ws.Range("E25").Value = "  -0.57%  "
$ws.Range("E26").Value = "  -1.88%  "
$ws.Range("D27").Value = "'1.66"
$ws.Range("E27").Value = "  -6.79%  "
$ws.Range("B28").Value = "Cosmos"
$ws.Range("C28").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D28").Value = "'6.46"
$ws.Range("E28").Value = "  -3.63%  "
$ws.Range("B29").Value = "EthereumClassic"
$ws.Range("C29").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D29").Value = "'14.86"
$ws.Range("E29").Value = "  -1.54%  "
$ws.Range("E30").Value = "  -6.01%  "
$ws.Range("D31").Value = "'0.0463"
$ws.Range("E31").Value = "  -3.53%  "
$ws.Range("E32").Value = "  -2.10%  "
$ws.Range("D33").Value = "'3.00"
$ws.Range("E33").Value = "  -3.29%  "
$ws.Range("E34").Value = "  -1.26%  "
$ws.Range("E35").Value = "  -4.06%  "
$ws.Range("D36").Value = "1.091.42"
$ws.Range("E36").Value = "  -3.49%  "
$ws.Range("E37").Value = "  -0.70%  "
$ws.Range("B38").Value = "VeChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D38").Value = "'0.0151"
$ws.Range("E38").Value = "  -1.73%  "
$ws.Range("B39").Value = "MXToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D39").Value = "'2.32"
$ws.Range("E39").Value = "  -4.65%  "
$ws.Range("D40").Value = "'0.497"
$ws.Range("E40").Value = "  -3.36%  "
$ws.Range("D41").Value = "'0.775"
$ws.Range("E41").Value = "  -8.62%  "
$ws.Range("D42").Value = "'0.766"
$ws.Range("E42").Value = "  -0.82%  "
$ws.Range("D43").Value = "'92.61"
$ws.Range("E43").Value = "  -5.56%  "
$ws.Range("E44").Value = "  -2.43%  "
$ws.Range("D45").Value = "1.701.02"
$ws.Range("E45").Value = "  -3.60%  "
$ws.Range("D46").Value = "0.0₆0109"
$ws.Range("E46").Value = "  -4.58%  "
$ws.Range("D47").Value = "'52.80"
$ws.Range("E47").Value = "  -3.19%  "
$ws.Range("E48").Value = "  -4.07%  "
$ws.Range("E49").Value = "  -3.64%  "
$ws.Range("E50").Value = "  -1.76%  "
$ws.Range("E51").Value = "  -0.56%  "
